$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.3795
$ws.Range("D3").Value = -8.713800000000001
$ws.Range("E3").Value = 16.25419999999999
$ws.Range("D4").Value = -7.405199999999999
$ws.Range("E9").Value = 17.35350000000001
$ws.Range("A11").Value = -21.80720000000001
$ws.Range("A12").Value = -21.54900000000001
$ws.Range("D14").Value = -7.433600000000005
$ws.Range("A15").Value = -21.71449999999999
$ws.Range("E15").Value = 16.35840000000001
$ws.Range("E19").Value = 16.20000000000001
$ws.Range("E20").Value = 16.00509999999999
$ws.Range("E25").Value = 17.17010000000001
$ws.Range("D26").Value = -8.744200000000006
$ws.Range("A27").Value = -21.9765
$ws.Range("E27").Value = 16.5029
$ws.Range("A28").Value = -22.02270000000001
$ws.Range("E28").Value = 16.51600000000001
$ws.Range("E30").Value = 15.7892
$ws.Range("A31").Value = -21.85750000000001
$ws.Range("D31").Value = -9.241199999999997
$ws.Range("A32").Value = -21.39919999999999
$ws.Range("E32").Value = 15.493
$ws.Range("D35").Value = -9.120699999999992
$ws.Range("A36").Value = -20.39780000000001
$ws.Range("D37").Value = -7.983599999999998
$ws.Range("A38").Value = -20.38980000000002
$ws.Range("D39").Value = -7.001500000000004
$ws.Range("D40").Value = -7.691199999999996
$ws.Range("E44").Value = 16.11119999999999
$ws.Range("D45").Value = -7.262799999999999
$ws.Range("A46").Value = -21.87449999999999
$ws.Range("E47").Value = 16.4942
$ws.Range("D52").Value = -7.592099999999997
$ws.Range("A54").Value = -21.58969999999999
$ws.Range("A55").Value = -22.29780000000001
$ws.Range("A56").Value = -22.12390000000001
$ws.Range("D57").Value = -8.619900000000003
$ws.Range("E58").Value = 16.9456
$ws.Range("E62").Value = 16.73659999999999
$ws.Range("A67").Value = -21.44609999999998
$ws.Range("A69").Value = -21.68569999999998
$ws.Range("A72").Value = -21.77810000000001
$ws.Range("A73").Value = -19.71629999999999
$ws.Range("E77").Value = 16.98180000000001
$ws.Range("E78").Value = 16.52700000000003
$ws.Range("D81").Value = -7.326600000000001
$ws.Range("A83").Value = -21.92519999999999
$ws.Range("D83").Value = -8.825699999999994
$ws.Range("E84").Value = 16.6605
$ws.Range("A86").Value = -22.05630000000001
$ws.Range("E89").Value = 17.35080000000002
$ws.Range("A91").Value = -21.43400000000002
$ws.Range("E91").Value = 17.84860000000002
$ws.Range("E92").Value = 17.97320000000002
$ws.Range("A93").Value = -21.13509999999999
$ws.Range("E96").Value = 15.81369999999999
$ws.Range("A99").Value = -19.8203
$ws.Range("D100").Value = -8.3856
$ws.Range("D102").Value = -7.730599999999998
$ws.Range("E102").Value = 16.5444
